$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values for rows 2-8, columns B-G, per regenerated s_vals data
# (filtered save games)

$ws.Range("B2").Value = 3.230985683306322
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.553084769722144

$ws.Range("B3").Value = 1.459612070389937
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 0.1575252929769615
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 3.781711156805759

$ws.Range("B4").Value = 1.459612070389937
$ws.Range("C4").Value = 1.667794583268128
$ws.Range("D4").Value = 0.8054896365839992
$ws.Range("E4").Value = 8.660232485948974
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.59312877619104

$ws.Range("B5").Value = 0.01514828764759746
$ws.Range("C5").Value = 0.04240448674262143
$ws.Range("D5").Value = 0.1575252929769615
$ws.Range("E5").Value = 0.496779210170732
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0.7118572775379124

$ws.Range("B6").Value = 3.230985683306322
$ws.Range("C6").Value = 1.667794583268128
$ws.Range("D6").Value = 26.21740644021617
$ws.Range("E6").Value = 0.496779210170732
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 31.61296591696135

$ws.Range("B7").Value = 3.230985683306322
$ws.Range("C7").Value = 1.667794583268128
$ws.Range("D7").Value = 3.900430680208489
$ws.Range("E7").Value = 0.496779210170732
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 9.295990156953671

$ws.Range("B8").Value = 3.230985683306322
$ws.Range("C8").Value = 1.667794583268128
$ws.Range("D8").Value = 3.900430680208489
$ws.Range("E8").Value = 0.496779210170732
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 9.295990156953671
